# Applies the "before -> after" edit described by the diff:
#  1) Several data rows (identified by the "Beteckning" in column A) were
#     reordered within the sheet - i.e. whole rows swapped places with
#     each other. We implement this as a row-content permutation: for
#     each affected destination row we snapshot the ENTIRE row (A:Z,
#     values AND formulas) from its current ("before") location first,
#     and only afterwards write the snapshots into their new homes - so
#     a row is never overwritten before it has been captured.
#  2) Every data row's "Förändrad" value (column C) advances by one day
#     (46076 -> 46077), independent of whether that row moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: snapshot the full A:Z content (values + formulas) of every
# row that ends up moving, BEFORE any write happens.
# ---------------------------------------------------------------------
$snapRow6  = $ws.Range("A6:Z6").Formula
$snapRow5  = $ws.Range("A5:Z5").Formula

$snapRow12 = $ws.Range("A12:Z12").Formula
$snapRow11 = $ws.Range("A11:Z11").Formula

$snapRow14 = $ws.Range("A14:Z14").Formula
$snapRow13 = $ws.Range("A13:Z13").Formula

$snapRow21 = $ws.Range("A21:Z21").Formula
$snapRow22 = $ws.Range("A22:Z22").Formula
$snapRow23 = $ws.Range("A23:Z23").Formula

# Rows 31-37 never carry anything past column Q (no species list / no
# links for these particular findings), so restrict the snapshot to
# A:Q - this avoids gratuitously touching the already-blank R:Z cells
# (re-assigning an empty value to them would needlessly flip their
# underlying storage from an empty shared string to a truly-empty
# cell).
$snapRow31 = $ws.Range("A31:Q31").Formula
$snapRow32 = $ws.Range("A32:Q32").Formula
$snapRow33 = $ws.Range("A33:Q33").Formula
$snapRow35 = $ws.Range("A35:Q35").Formula
$snapRow36 = $ws.Range("A36:Q36").Formula
$snapRow37 = $ws.Range("A37:Q37").Formula

# ---------------------------------------------------------------------
# Step 2: write the snapshots back into their new row positions.
# ---------------------------------------------------------------------

# Rows 5 <-> 6  (A 21219-2023 <-> A 27865-2024)
$ws.Range("A5:Z5").Formula  = $snapRow6
$ws.Range("A6:Z6").Formula  = $snapRow5

# Rows 11 <-> 12 (A 38039-2022 <-> A 27636-2023)
$ws.Range("A11:Z11").Formula = $snapRow12
$ws.Range("A12:Z12").Formula = $snapRow11

# Rows 13 <-> 14 (A 43714-2025 <-> A 39320-2024)
$ws.Range("A13:Z13").Formula = $snapRow14
$ws.Range("A14:Z14").Formula = $snapRow13

# Rows 21, 22, 23 rotate: new21=old23, new22=old21, new23=old22
$ws.Range("A21:Z21").Formula = $snapRow23
$ws.Range("A22:Z22").Formula = $snapRow21
$ws.Range("A23:Z23").Formula = $snapRow22

# Rows 31..37 (34 stays fixed) rotate:
# new31=old33, new32=old31, new33=old37, new35=old32, new36=old35, new37=old36
$ws.Range("A31:Q31").Formula = $snapRow33
$ws.Range("A32:Q32").Formula = $snapRow31
$ws.Range("A33:Q33").Formula = $snapRow37
$ws.Range("A35:Q35").Formula = $snapRow32
$ws.Range("A36:Q36").Formula = $snapRow35
$ws.Range("A37:Q37").Formula = $snapRow36

# ---------------------------------------------------------------------
# Step 3: bump the "Förändrad" date (column C) for every data row
# (2-37) from 46076 to 46077, regardless of whether the row moved.
# ---------------------------------------------------------------------
$ws.Range("C2:C37").Value2 = 46077

# ---------------------------------------------------------------------
# Step 4: re-assigning wrapped, multi-line text (column R) can make the
# host recompute "best fit" row heights. The source rows all used a
# fixed explicit height of 15 points, independent of how much text
# wraps in column R, so restore that explicit height on every row
# whose content we rewrote above (only rows that carry column-R text
# are actually affected, but it is harmless to reset all of them).
# ---------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight  = 15
$ws.Rows.Item(6).RowHeight  = 15
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 15
$ws.Rows.Item(21).RowHeight = 15
$ws.Rows.Item(22).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 15
